$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 676, shifting existing rows 676-717 down to 677-718
$ws.Rows("676:676").Insert()

# Populate the newly inserted row with the new data point.
# The leading apostrophe forces the date-like text to be stored as literal
# text instead of being auto-converted into a date serial number.
$ws.Range("A676").Value = "'2026/01/18"
$ws.Range("B676").Value = "日"
$ws.Range("C676").Value = 19
$ws.Range("D676").Value = 166
